# The deck currently ships two theme parts:
#   theme1.xml -> "Office Theme" / clrScheme "Office"      (orphaned; only
#                  referenced by the notes master, not rendered for slides)
#   theme2.xml -> "Integral"     / clrScheme "Red Violet"  (the live theme,
#                  referenced by the slide master that all slides use)
#
# The target edit swaps the two themes' content, so the slides' live theme
# becomes the plain "Office" colour scheme instead of "Red Violet".
#
# Re-point the presentation's (single, shared) theme colour scheme -
# reachable via $p.SlideMaster.Theme.ThemeColorScheme - to the "Office"
# palette that used to live in the unused theme part.

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.Theme.ThemeColorScheme

function Hex2RGB($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# ThemeColorScheme.Item index order: 1=dk1 2=lt1 3=dk2 4=lt2 5=accent1
# 6=accent2 7=accent3 8=accent4 9=accent5 10=accent6 11=hlink 12=folHlink
$officeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

for ($i = 0; $i -lt $officeColors.Length; $i++) {
    $cs.Item($i + 1).RGB = Hex2RGB $officeColors[$i]
}
